# license-master.pptx - "Initial prep for SC21" edit
#
# 1. Update the requested-citation paragraph (slide 1, Content Placeholder 2,
#    3rd paragraph) to the SC21 citation + new DOI.
# 2. Update the cached "datetimeFigureOut" footer date shown on the
#    Handout Master / Notes Master from 8/7/2021 to 9/1/2021.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- 1. Citation text on the License/Citation slide -----------------------
$contentShape = $s.Shapes.Item(2)
$citationPara = $contentShape.TextFrame.TextRange.Paragraphs(3)

$citationPara.Runs(1).Text = "The requested citation the overall tutorial is: David E. Bernholdt, Anshu Dubey, Patricia A. Grubel, Rinku K. Gupta, and Gregory R. Watson, Better Scientific Software tutorial, in the International Conference for High-Performance Computing, Networking, Storage, and Analysis (SC21), St. Louis, MO, USA and online, 2021. DOI: "
$citationPara.Runs(2).Text = "10.6084/m9.figshare.16556628"

# --- 2. Footer date on Handout Master / Notes Master -----------------------
$p.HandoutMaster.HeadersFooters.DateAndTime.Text = "9/1/2021"
$p.NotesMaster.HeadersFooters.DateAndTime.Text = "9/1/2021"
